# Update column F (dSF) values for rows 2-37 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -2
    3  = 3
    4  = -8
    5  = -2
    7  = -2
    9  = -2
    10 = 2
    11 = 1
    12 = 5
    13 = 2
    14 = 5
    15 = -3
    16 = -4
    17 = 3
    18 = -1
    19 = 1
    20 = 2
    21 = -4
    22 = 4
    23 = -3
    24 = -2
    25 = -1
    26 = 0
    27 = -2
    29 = -2
    30 = -4
    31 = 2
    32 = 7
    33 = 1
    34 = -1
    35 = 6
    36 = -1
    37 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
